$d = $word.ActiveDocument

# Locate the empty paragraph that immediately follows the "Commit 10" body
# text ("...style.less file") -- this is the insertion anchor. Two new
# paragraphs (a bold/underlined "Commit 10:" heading, plus a follow-up body
# paragraph) must land directly before it, while the anchor paragraph
# itself is left untouched immediately afterwards.
$anchorIndex = 0
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "`r") {
        $prevText = $d.Paragraphs.Item($i - 1).Range.Text
        if ($prevText -like "*style.less file*") {
            $anchorIndex = $i
            break
        }
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex)

# NOTE: InsertBefore always lands text immediately ahead of the anchor
# paragraph's current position, so inserting twice through the same
# paragraph reference would reverse the intended order if done
# heading-then-body. Insert the body text first, then the heading, so
# the heading ends up first in the final document.

$quote1 = [char]0x201C
$quote2 = [char]0x201D
$bodyText = "This step is still left we have to do this also.. Please add this add this new folder with this navigation.js file in master branch " + $quote1 + "src/assets/js/ navigation.js" + $quote2
$anchor.Range.InsertBefore($bodyText + "`r")
$anchor.Range.InsertBefore("Commit 10: Missing 1 no commit mobile nav issue`r")

# The heading is now directly at $anchorIndex, the body right after it,
# and the untouched original anchor paragraph follows that.
$headingIndex = $anchorIndex
$bodyIndex = $anchorIndex + 1

$headingRange = $d.Paragraphs.Item($headingIndex).Range
$headingRange.Font.Bold = $true
$headingRange.Font.Size = 14
$headingRange.Font.Underline = 1

$bodyPara = $d.Paragraphs.Item($bodyIndex)
$bodyRange = $bodyPara.Range
$bodyRange.Font.Size = 14
$bodyPara.Range.ParagraphFormat.FirstLineIndent = 36

# Apply yellow highlight to "js/", the following space, and "navigation.js"
# inside the body paragraph, matching the original formatting.
$bodyStart = $bodyPara.Range.Start
$prefix = "This step is still left we have to do this also.. Please add this add this new folder with this navigation.js file in master branch " + $quote1 + "src/assets/"
$jsStart = $bodyStart + $prefix.Length
$jsEnd = $jsStart + ("js/").Length
$d.Range($jsStart, $jsEnd).HighlightColorIndex = 7

$spaceEnd = $jsEnd + 1
$spaceRange = $d.Range($jsEnd, $spaceEnd)
$spaceRange.HighlightColorIndex = 7
$spaceRange.Font.Size = 10

$navStart = $spaceEnd
$navEnd = $navStart + ("navigation.js").Length
$d.Range($navStart, $navEnd).HighlightColorIndex = 7

Write-Output "ok"
